# Fix table DCI export bug: correct the sizing-worksheet inputs that feed
# the "One Table per Node" calculations (column E of the table block plus
# the DCI/node/retention inputs in column C further down).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("One Table per Node")

# Number of int columns in the table (row 4)
$ws.Range("E4").Value = 2

# Max varchar size (row 5)
$ws.Range("E5").Value = 510

# Number of nullable columns (row 7)
$ws.Range("E7").Value = 7

# Number of DCIs (row 12)
$ws.Range("C12").Value = 12000

# Number of nodes (row 13)
$ws.Range("C13").Value = 400

# Retention, in days (row 15)
$ws.Range("C15").Value = 90
